# Updated cryptos list with latest prices / 1h volume change figures
# (GitHub Actions scheduled refresh). Two coin pairs (USDC/LidoStakedEther,
# Uniswap/WrappedBTC, Litecoin/InternetComputer(DFINITY),
# Stacks/InjectiveProtocol, VeChain/Fetch.AI) swapped rank order as their
# prices moved past each other.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.569.18"
$ws.Range("E2").Value = "  -4.91%  "

# Row 3
$ws.Range("D3").Value = "'3.323.76"
$ws.Range("E3").Value = "  -5.95%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'549.99"
$ws.Range("E5").Value = "  -2.85%  "

# Row 6
$ws.Range("D6").Value = "'170.71"
$ws.Range("E6").Value = "  -8.87%  "

# Row 7
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  -4.89%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "'3.316.87"
$ws.Range("E8").Value = "  -6.09%  "

# Row 9
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").Value = "'0.612"
$ws.Range("E10").Value = "  -5.36%  "

# Row 11
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -4.76%  "

# Row 12
$ws.Range("D12").Value = "'53.05"
$ws.Range("E12").Value = "  -5.52%  "

# Row 13
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  -6.48%  "

# Row 14
$ws.Range("E14").Value = "  -6.56%  "

# Row 15
$ws.Range("D15").Value = "'3.845.74"
$ws.Range("E15").Value = "  -5.95%  "

# Row 16
$ws.Range("E16").Value = "  -4.18%  "

# Row 17
$ws.Range("D17").Value = "'17.77"
$ws.Range("E17").Value = "  -5.84%  "

# Row 18
$ws.Range("D18").Value = "'3.318.31"
$ws.Range("E18").Value = "  -5.90%  "

# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'63.507.80"
$ws.Range("E19").Value = "  -5.03%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'11.58"
$ws.Range("E20").Value = "  -5.08%  "

# Row 21
$ws.Range("D21").Value = "'0.964"
$ws.Range("E21").Value = "  -4.64%  "

# Row 22
$ws.Range("D22").Value = "'405.48"
$ws.Range("E22").Value = "  -4.14%  "

# Row 23
$ws.Range("E23").Value = "  -2.07%  "

# Row 24
$ws.Range("D24").Value = "'4.26"
$ws.Range("E24").Value = "  -1.11%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.20"
$ws.Range("E25").Value = "  +5.59%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'82.53"
$ws.Range("E26").Value = "  -4.72%  "

# Row 27
$ws.Range("D27").Value = "'10.55"
$ws.Range("E27").Value = "  -4.86%  "

# Row 28
$ws.Range("E28").Value = "  -8.43%  "

# Row 29
$ws.Range("D29").Value = "'8.57"
$ws.Range("E29").Value = "  -6.99%  "

# Row 30
$ws.Range("D30").Value = "'28.97"
$ws.Range("E30").Value = "  -5.47%  "

# Row 31
$ws.Range("D31").Value = "'6.38"
$ws.Range("E31").Value = "  -6.01%  "

# Row 32
$ws.Range("D32").Value = "'11.28"
$ws.Range("E32").Value = "  -5.63%  "

# Row 33
$ws.Range("D33").Value = "'571.65"
$ws.Range("E33").Value = "  -8.65%  "

# Row 34
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -6.39%  "

# Row 35
$ws.Range("D35").Value = "'57.17"
$ws.Range("E35").Value = "  -5.42%  "

# Row 36
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("E37").Value = "  -3.79%  "

# Row 38
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'35.09"
$ws.Range("E38").Value = "  -9.21%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
$ws.Range("D40").Value = "'0.0₃0735"
$ws.Range("E40").Value = "  -11.11%  "

# Row 41
$ws.Range("D41").Value = "'3.156.82"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("D42").Value = "'0.365"
$ws.Range("E42").Value = "  -6.57%  "

# Row 43
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -2.63%  "

# Row 45
$ws.Range("D45").Value = "'3.18"
$ws.Range("E45").Value = "  -4.52%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.43"
$ws.Range("E46").Value = "  -8.49%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0401"
$ws.Range("E47").Value = "  -5.01%  "

# Row 48
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -5.14%  "

# Row 49
$ws.Range("D49").Value = "'0.128"
$ws.Range("E49").Value = "  -4.92%  "

# Row 50
$ws.Range("D50").Value = "'132.77"
$ws.Range("E50").Value = "  -5.12%  "

# Row 51
$ws.Range("D51").Value = "'7.99"
$ws.Range("E51").Value = "  -7.28%  "
